$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update A8: "Culebra" -> "yami"
$ws.Range("A8").Value = "yami"

# Highlight header row B3:E3 with a yellow fill
$ws.Range("B3:E3").Interior.Color = 65535

# Update selection to C6
$ws.Range("C6").Select()
